$wb = $excel.ActiveWorkbook

# Use the existing header formatting (bold, border, centered) from sheet 1
# as the template for the new sheets' header rows, so we don't introduce any
# new style/font entries into styles.xml.
$headerTemplate = $wb.Worksheets.Item("PAGE COUNT (DATOPIAN)").Range("A1:B1")

# Insert the new sheets right after the existing "PAGE COUNT (AIR)" sheet,
# preserving order: ... , PAGE COUNT (AIR), RESOURCE COUNT (DATOPIAN ONLY), RESOURCE COUNT (AIR ONLY)
$pageCountAir = $wb.Worksheets.Item("PAGE COUNT (AIR)")

# --- Sheet: RESOURCE COUNT (DATOPIAN ONLY) ---
$wsDatopian = $wb.Worksheets.Add($null, $pageCountAir)
$wsDatopian.Name = "RESOURCE COUNT (DATOPIAN ONLY)"

$datopianData = @(
    @("domain", "resource count"),
    @("stateaid.nysed.gov", 101),
    @("data.nysed.gov", 39),
    @("results.ed.gov", 25),
    @("eric.ed.gov", 22),
    @("p12.nysed.gov", 14),
    @("emsc32.nysed.gov", 8),
    @("pdp.ed.gov", 6),
    @("americanenglish.state.gov", 2),
    @("portal.nysed.gov", 2),
    @("ccsso.org", 1),
    @("p1232.nysed.gov", 1)
)

for ($i = 0; $i -lt $datopianData.Count; $i++) {
    $row = $i + 1
    $wsDatopian.Cells.Item($row, 1).Value = $datopianData[$i][0]
    $wsDatopian.Cells.Item($row, 2).Value = $datopianData[$i][1]
}

$headerTemplate.Copy()
$wsDatopian.Range("A1:B1").PasteSpecial(-4122)

# --- Sheet: RESOURCE COUNT (AIR ONLY) ---
$wsAir = $wb.Worksheets.Add($null, $wsDatopian)
$wsAir.Name = "RESOURCE COUNT (AIR ONLY)"

$airData = @(
    @("domain", "resource count"),
    @("fsadownload.ed.gov", 215),
    @("fp.ed.gov", 59),
    @("charterschoolcenter.ed.gov", 32),
    @("y4y.ed.gov", 30),
    @("surveys.ope.ed.gov", 11),
    @("nslds.ed.gov", 7),
    @("experimentalsites.ed.gov", 1),
    @("nsldsfap.ed.gov", 1),
    @("ope.ed.gov", 1)
)

for ($i = 0; $i -lt $airData.Count; $i++) {
    $row = $i + 1
    $wsAir.Cells.Item($row, 1).Value = $airData[$i][0]
    $wsAir.Cells.Item($row, 2).Value = $airData[$i][1]
}

$headerTemplate.Copy()
$wsAir.Range("A1:B1").PasteSpecial(-4122)

# Restore the first sheet as active/selected (matches the original workbook state)
$wb.Worksheets.Item("PAGE COUNT (DATOPIAN)").Activate()
